$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 388.75
$ws.Range("I2").Value = 333.2857
$ws.Range("K2").Value = 333.2857
$ws.Range("M2").Value = -220.2857
$ws.Range("H9").Value = 85
$ws.Range("J9").Value = 93.25
$ws.Range("L9").Value = 93.25
$ws.Range("N9").Value = -431.25
$ws.Range("H17").Value = 975.8570999999999
$ws.Range("I17").Value = 377.75
$ws.Range("K17").Value = 1133.25
$ws.Range("M17").Value = -965.25
$ws.Range("H53").Value = 476.08
$ws.Range("I53").Value = 363.75
$ws.Range("K53").Value = 363.75
$ws.Range("M53").Value = 273.25
$ws.Range("H80").Value = 805.2857
$ws.Range("I80").Value = 287.2
$ws.Range("K80").Value = 861.5999999999999
$ws.Range("M80").Value = 136.4000000000001
$ws.Range("H83").Value = 805.2857
$ws.Range("I83").Value = 287.2
$ws.Range("K83").Value = 2584.8
$ws.Range("M83").Value = 2407.2
$ws.Range("H133").Value = 59995
$ws.Range("J133").Value = 59995
$ws.Range("L133").Value = 59995
$ws.Range("N133").Value = -70115
$ws.Range("H141").Value = 7830.7144
$ws.Range("I141").Value = 5764.857
$ws.Range("K141").Value = 17294.571
$ws.Range("M141").Value = -12114.571
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2078.5518
$ws.Range("I2").Value = 1676.7059
$ws.Range("J2").Value = 2647.8333
$ws.Range("K2").Value = 1676.7059
$ws.Range("L2").Value = 2647.8333
$ws.Range("M2").Value = -1563.7059
$ws.Range("N2").Value = -2873.8333
$ws.Range("H32").Value = 2540.8462
$ws.Range("I32").Value = 2643.8044
$ws.Range("K32").Value = 2643.8044
$ws.Range("M32").Value = -2356.8044
$ws.Range("H61").Value = 6086.3335
$ws.Range("I61").Value = 6086.3335
$ws.Range("K61").Value = 6086.3335
$ws.Range("M61").Value = -5874.3335
$ws.Range("H116").Value = 2078.5518
$ws.Range("I116").Value = 1676.7059
$ws.Range("J116").Value = 2647.8333
$ws.Range("K116").Value = 1676.7059
$ws.Range("L116").Value = 2647.8333
$ws.Range("M116").Value = 617.2941000000001
$ws.Range("N116").Value = -7235.8333
$ws.Range("H132").Value = 2190.077
$ws.Range("I132").Value = 2105.348
$ws.Range("K132").Value = 6316.044
$ws.Range("M132").Value = -3786.044
$ws.Range("H136").Value = 6086.3335
$ws.Range("I136").Value = 6086.3335
$ws.Range("K136").Value = 18259.0005
$ws.Range("M136").Value = -15709.0005
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2078.5518
$ws.Range("I3").Value = 1676.7059
$ws.Range("J3").Value = 2647.8333
$ws.Range("K3").Value = 1676.7059
$ws.Range("L3").Value = 2647.8333
$ws.Range("M3").Value = -1562.7059
$ws.Range("N3").Value = -2875.8333
$ws.Range("H22").Value = 1899.091
$ws.Range("I22").Value = 1887.5555
$ws.Range("K22").Value = 1887.5555
$ws.Range("M22").Value = -1714.5555
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("M26").Value = ""
$ws.Range("H86").Value = 21741686
$ws.Range("I86").Value = 2322.4443
$ws.Range("J86").Value = 100003400
$ws.Range("K86").Value = 2322.4443
$ws.Range("L86").Value = 100003400
$ws.Range("M86").Value = -1199.4443
$ws.Range("N86").Value = -100005646
$ws.Range("H89").Value = 21741686
$ws.Range("I89").Value = 2322.4443
$ws.Range("J89").Value = 100003400
$ws.Range("K89").Value = 11612.2215
$ws.Range("L89").Value = 500017000
$ws.Range("M89").Value = -5996.2215
$ws.Range("N89").Value = -500028232
$ws.Range("H100").Value = 75000
$ws.Range("J100").Value = 75000
$ws.Range("L100").Value = 75000
$ws.Range("N100").Value = -77164
$ws.Range("H132").Value = 56832
$ws.Range("J132").Value = 56832
$ws.Range("L132").Value = 56832
$ws.Range("N132").Value = -66952
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3814.647
$ws.Range("I132").Value = 4414.5
$ws.Range("J132").Value = 2957.7144
$ws.Range("K132").Value = 13243.5
$ws.Range("L132").Value = 8873.143199999999
$ws.Range("M132").Value = -10713.5
$ws.Range("N132").Value = -13933.1432
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 3913.375
$ws.Range("I3").Value = 3913.375
$ws.Range("K3").Value = 11740.125
$ws.Range("M3").Value = -11628.125
$ws.Range("H5").Value = 1280
$ws.Range("I5").Value = 713.53845
$ws.Range("K5").Value = 2140.61535
$ws.Range("M5").Value = -2028.61535
$ws.Range("H107").Value = 1660.7826
$ws.Range("I107").Value = 4108.7144
$ws.Range("J107").Value = 589.8125
$ws.Range("K107").Value = 12326.1432
$ws.Range("L107").Value = 1769.4375
$ws.Range("M107").Value = -10406.1432
$ws.Range("N107").Value = -5609.4375
$ws.Range("H129").Value = 3256.7144
$ws.Range("J129").Value = 4225.4
$ws.Range("L129").Value = 12676.2
$ws.Range("N129").Value = -22676.2
$ws.Range("H135").Value = 1280
$ws.Range("I135").Value = 713.53845
$ws.Range("K135").Value = 6421.84605
$ws.Range("M135").Value = -3886.84605
$ws.Range("H140").Value = 1751
$ws.Range("I140").Value = 1411.125
$ws.Range("K140").Value = 4233.375
$ws.Range("M140").Value = 946.625
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H64").Value = 62500
$ws.Range("J64").Value = 62500
$ws.Range("L64").Value = 62500
$ws.Range("N64").Value = -62996
$ws.Range("H67").Value = 62500
$ws.Range("J67").Value = 62500
$ws.Range("L67").Value = 62500
$ws.Range("N67").Value = -64216
$ws.Range("H97").Value = 713.55
$ws.Range("I97").Value = 611.9231
$ws.Range("J97").Value = 902.2857
$ws.Range("K97").Value = 611.9231
$ws.Range("L97").Value = 902.2857
$ws.Range("M97").Value = -115.9231
$ws.Range("N97").Value = -1894.2857
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = ""
$ws.Range("N104").Value = 0
$ws.Range("H122").Value = 6349.1
$ws.Range("J122").Value = 5014.143
$ws.Range("L122").Value = 15042.429
$ws.Range("N122").Value = -19942.429
$ws.Range("H123").Value = 69999.25
$ws.Range("J123").Value = 69999.25
$ws.Range("L123").Value = 69999.25
$ws.Range("N123").Value = -74899.25
$ws.Range("H126").Value = 2447.5
$ws.Range("I126").Value = 2447.5
$ws.Range("K126").Value = 7342.5
$ws.Range("M126").Value = -4872.5
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3607.6086
$ws.Range("I22").Value = 3568.375
$ws.Range("K22").Value = 3568.375
$ws.Range("M22").Value = -3273.375
$ws.Range("H27").Value = 3607.6086
$ws.Range("I27").Value = 3568.375
$ws.Range("K27").Value = 3568.375
$ws.Range("M27").Value = -3461.375
$ws.Range("H46").Value = 1915.1538
$ws.Range("I46").Value = 986.5
$ws.Range("J46").Value = 2084
$ws.Range("K46").Value = 986.5
$ws.Range("L46").Value = 2084
$ws.Range("M46").Value = -798.5
$ws.Range("N46").Value = -2460
$ws.Range("H68").Value = 12129.1
$ws.Range("J68").Value = 14598.75
$ws.Range("L68").Value = 14598.75
$ws.Range("N68").Value = -16096.75
$ws.Range("H71").Value = 12129.1
$ws.Range("J71").Value = 14598.75
$ws.Range("L71").Value = 72993.75
$ws.Range("N71").Value = -80481.75
$ws.Range("H93").Value = 3240.2896
$ws.Range("I93").Value = 1701.9546
$ws.Range("K93").Value = 1701.9546
$ws.Range("M93").Value = -453.9546
$ws.Range("H98").Value = 45355
$ws.Range("J98").Value = 45355
$ws.Range("L98").Value = 45355
$ws.Range("N98").Value = -51345
$ws.Range("H132").Value = 4208.722
$ws.Range("I132").Value = 4037.3845
$ws.Range("J132").Value = 4367.8213
$ws.Range("K132").Value = 12112.1535
$ws.Range("L132").Value = 13103.4639
$ws.Range("M132").Value = -9582.1535
$ws.Range("N132").Value = -18163.4639
$ws.Range("H133").Value = 29999
$ws.Range("J133").Value = 29999
$ws.Range("L133").Value = 29999
$ws.Range("N133").Value = -35059
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1198.5
$ws.Range("I81").Value = 1138.2
$ws.Range("J81").Value = 1500
$ws.Range("K81").Value = 2276.4
$ws.Range("L81").Value = 3000
$ws.Range("M81").Value = -1215.4
$ws.Range("N81").Value = -5122
$ws.Range("H84").Value = 1198.5
$ws.Range("I84").Value = 1138.2
$ws.Range("J84").Value = 1500
$ws.Range("K84").Value = 11382
$ws.Range("L84").Value = 15000
$ws.Range("M84").Value = -6078
$ws.Range("N84").Value = -25608
$ws.Range("H132").Value = 7412.92
$ws.Range("I132").Value = 5339
$ws.Range("K132").Value = 16017
$ws.Range("M132").Value = -13487
